# Added seasonal statistics and a combined RSTs value sheet
# (done in a very ugly way, refactor later)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title cell: "1979-2016" -> "1979-2016 00Z" -----------------------
$ws.Range("A1").Value = "1979-2016 00Z"

# --- Swap the West/Central row labels in each of the three blocks -----
# NCEP block (rows 2-5)
$ws.Range("A4").Value = "NCEP-Central"
$ws.Range("A5").Value = "NCEP-West"

# ERA block (rows 7-10)
$ws.Range("A9").Value = "ERA-Central"
$ws.Range("A10").Value = "ERA-West"

# ERA 2.5 block (rows 12-15)
$ws.Range("A14").Value = "ERA 2.5 -Central"
$ws.Range("A15").Value = "ERA 2.5 -West"

# --- Fix subtotal for the NCEP block (row 6) to skip the header row ---
$ws.Range("U6").Formula = "=SUM(U3:U5)"

# --- Add the per-row totals (column U) that were missing for the ERA
#     and ERA 2.5 blocks -------------------------------------------------
$ws.Range("U7").Formula = "=SUM(B7:T7)"
$ws.Range("U8").Formula = "=SUM(B8:T8)"
$ws.Range("U9").Formula = "=SUM(B9:T9)"
$ws.Range("U10").Formula = "=SUM(B10:T10)"

$ws.Range("U12").Formula = "=SUM(B12:T12)"
$ws.Range("U13").Formula = "=SUM(B13:T13)"
$ws.Range("U14").Formula = "=SUM(B14:T14)"
$ws.Range("U15").Formula = "=SUM(B15:T15)"

# --- Add the block-subtotal rows for ERA (row 11) and ERA 2.5 (row 16) -
$ws.Range("U11").Formula = "=SUM(U8:U10)"
$ws.Range("U16").Formula = "=SUM(U13:U15)"

# --- View state: zoom to 150% and select the new last row --------------
$excel.ActiveWindow.Zoom = 150
$ws.Range("A16").Select()
